$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 3 new rows at row 11 (table grows from 17 to 20 data+header rows)
# ---------------------------------------------------------------------------
$ws.Rows("11:13").Insert()

# The 3-row insert shifted the old rows 11-17 down to rows 14-20, which means
# some cells that used to hold a second-author (column B) value now sit in
# rows that must end up single-author. Wipe the whole block first so every
# cell below gets written from a clean slate.
$ws.Range("A11:J20").ClearContents()

# ---------------------------------------------------------------------------
# 2) Write the final (post-sort) contents for rows 11-20. The three brand new
#    records ("Faure F" / "Medel C" x2) are merged alphabetically with the
#    records that used to occupy rows 11-17, so most of this block is simply
#    re-writing pre-existing rows that shifted position.
# ---------------------------------------------------------------------------

# Row 11 - new record: Faure F / Medel C
$ws.Range("A11").Value = "Faure F"
$ws.Range("B11").Value = "Medel C"
$ws.Range("E11").Value = 2020
$ws.Range("F11").Value = "Does the Exposure to the Business Cycle Improve Consumer Perceptions for Forecasting? Microdata Evidence from Chile"
$ws.Range("G11").Value = "Working Paper"
$ws.Range("H11").Value = "Documentos de Trabajo (Banco Central)"
$ws.Range("I11").Value = "Microeconomía"
$ws.Range("J11").Value = "https://www.bcentral.cl/en/content/-/details/working-papers-n-888"

# Row 12 - Figueroa C / Pedersen M (previously row 11)
$ws.Range("A12").Value = "Figueroa C"
$ws.Range("B12").Value = "Pedersen M"
$ws.Range("E12").Value = 2019
$ws.Range("F12").Value = "Extracting information on economic activity from business and consumer surveys in an emerging economy (Chile)"
$ws.Range("G12").Value = "Paper"
$ws.Range("H12").Value = "Economía Chilena"
$ws.Range("I12").Value = "Macroeconomía"
$ws.Range("J12").Value = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4041"

# Row 13 - García P / Olea S (previously row 12)
$ws.Range("A13").Value = "García P"
$ws.Range("B13").Value = "Olea S"
$ws.Range("E13").Value = 2015
$ws.Range("F13").Value = "Inversión Minera y Ajuste Macroeconómico en Australia y Chile"
$ws.Range("G13").Value = "Paper"
$ws.Range("H13").Value = "Economic Policy Papers (Banco Central)"
$ws.Range("I13").Value = "Macroeconomía"
$ws.Range("J13").Value = "https://ideas.repec.org/p/chb/bcchep/56.html"

# Row 14 - new record: Medel C (2018)
$ws.Range("A14").Value = "Medel C"
$ws.Range("E14").Value = 2018
$ws.Range("F14").Value = "Econometric Analysis on Survey-data-based Anchoring of Inflation Expectations in Chile"
$ws.Range("G14").Value = "Working Paper"
$ws.Range("H14").Value = "Documentos de Trabajo (Banco Central)"
$ws.Range("I14").Value = "Macroeconomía"
$ws.Range("J14").Value = "https://www.bcentral.cl/en/content/-/details/working-papers-n-825"

# Row 15 - new record: Medel C (2021)
$ws.Range("A15").Value = "Medel C"
$ws.Range("E15").Value = 2021
$ws.Range("F15").Value = "Searching for the Best Inflation Forecasters within a Consumer Perceptions Survey: Microdata Evidence from Chile"
$ws.Range("G15").Value = "Working Paper"
$ws.Range("H15").Value = "Documentos de Trabajo (Banco Central)"
$ws.Range("I15").Value = "Microeconomía"
$ws.Range("J15").Value = "https://www.bcentral.cl/en/content/-/details/working-papers-n-899"

# Row 16 - Pedersen M 2009 (previously row 13)
$ws.Range("A16").Value = "Pedersen M"
$ws.Range("E16").Value = 2009
$ws.Range("F16").Value = "Un indicador líder compuesto para la actividad económica en Chile"
$ws.Range("G16").Value = "Paper"
$ws.Range("H16").Value = "Monetaria"
$ws.Range("I16").Value = "Macroeconomía"
$ws.Range("J16").Value = "https://ideas.repec.org/a/cml/moneta/vxxxiiy2009i2p181-208.html"

# Row 17 - Pedersen M 2009 (previously row 14)
$ws.Range("A17").Value = "Pedersen M"
$ws.Range("E17").Value = 2009
$ws.Range("F17").Value = "Use of Chilean Business Surveys in Conjunctural Assessment and Short-term Forecasting"
$ws.Range("G17").Value = "Paper"
$ws.Range("H17").Value = "OECD workshop"
$ws.Range("I17").Value = "Macroeconomía"
$ws.Range("J17").Value = "https://www.google.com/url?sa=t&rct=j&q=&esrc=s&source=web&cd=&ved=2ahUKEwiE7Y60i5TvAhVQErkGHfecC4MQFjABegQIAhAD&url=http%3A%2F%2Fwww.oecd.org%2Fstd%2Fleading-indicators%2F43815334.pdf&usg=AOvVaw3BstLuhLtAOtjJeL5SsMj4"

# Row 18 - Pedersen M 2019 (previously row 15)
$ws.Range("A18").Value = "Pedersen M"
$ws.Range("E18").Value = 2019
$ws.Range("F18").Value = "Anomalies in macroeconomic prediction errors–evidence from Chilean private forecasters"
$ws.Range("G18").Value = "Paper"
$ws.Range("H18").Value = "International Journal of Forecasting"
$ws.Range("I18").Value = "Macroeconomía"
$ws.Range("J18").Value = "https://www.sciencedirect.com/science/article/abs/pii/S0169207019300676"

# Row 19 - Pincheira P (previously row 16)
$ws.Range("A19").Value = "Pincheira P"
$ws.Range("E19").Value = 2014
$ws.Range("F19").Value = "Predicción del Empleo Sectorial y Total en Base a Indicadores de Confianza Empresarial"
$ws.Range("G19").Value = "Paper"
$ws.Range("H19").Value = "Economía Chilena"
$ws.Range("I19").Value = "Macroeconomía"
$ws.Range("J19").Value = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/3564"

# Row 20 - Riquelme V / Riveros G (previously row 17)
$ws.Range("A20").Value = "Riquelme V"
$ws.Range("B20").Value = "Riveros G"
$ws.Range("E20").Value = 2018
$ws.Range("F20").Value = "Un Indicador Contemporáneo de Actividad (ICA) para Chile"
$ws.Range("G20").Value = "Paper"
$ws.Range("H20").Value = "Economía Chilena"
$ws.Range("I20").Value = "Macroeconomía"
$ws.Range("J20").Value = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4867"

# ---------------------------------------------------------------------------
# 3) Hyperlinks: the "Sitio web" column (J) carries a real hyperlink object
#    on top of the plain text already written above. Rebuild the hyperlink
#    list from scratch (rows 2-10 are untouched content-wise, but we still
#    need to keep them) so relationship ids come out in the same order as
#    the target file: J2..J9 unchanged, then J12,J13,J16..J20,J10,J15,J14,J11.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$links = [ordered]@{
    "J2"  = "https://ideas.repec.org/p/pra/mprapa/79809.html"
    "J3"  = "https://ideas.repec.org/p/pra/mprapa/83154.html"
    "J4"  = "https://doi.org/10.1080/02692171.2019.1645816"
    "J5"  = "https://www.bcentral.cl/en/content/-/details/monetary-policy-report-june-2015"
    "J6"  = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4043"
    "J7"  = "https://www.bcentral.cl/en/web/banco-central/content/-/detalle/documento-de-trabajo-n-883"
    "J8"  = "https://ideas.repec.org/a/chb/bcchni/v15y2012i1p105-117.html"
    "J9"  = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4042"
    "J12" = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4041"
    "J13" = "https://ideas.repec.org/p/chb/bcchep/56.html"
    "J16" = "https://ideas.repec.org/a/cml/moneta/vxxxiiy2009i2p181-208.html"
    "J17" = "https://www.google.com/url?sa=t&rct=j&q=&esrc=s&source=web&cd=&ved=2ahUKEwiE7Y60i5TvAhVQErkGHfecC4MQFjABegQIAhAD&url=http%3A%2F%2Fwww.oecd.org%2Fstd%2Fleading-indicators%2F43815334.pdf&usg=AOvVaw3BstLuhLtAOtjJeL5SsMj4"
    "J18" = "https://www.sciencedirect.com/science/article/abs/pii/S0169207019300676"
    "J19" = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/3564"
    "J20" = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4867"
    "J10" = "https://www.bcentral.cl/en/web/banco-central/content/-/detalle/analisis-de-sentimiento-basado-en-el-informe-de-percepciones-de-negocios-del-banco-central-de-chile"
    "J15" = "https://www.bcentral.cl/en/content/-/details/working-papers-n-899"
    "J14" = "https://www.bcentral.cl/en/content/-/details/working-papers-n-825"
    "J11" = "https://www.bcentral.cl/en/content/-/details/working-papers-n-888"
}

foreach ($ref in $links.Keys) {
    $ws.Hyperlinks.Add($ws.Range($ref), $links[$ref])
    $ws.Range($ref).Style = "Hipervínculo"
}

# ---------------------------------------------------------------------------
# 4) Restore the sort metadata (sortState) to cover the new A2:J20 range -
#    re-applying the same two-level sort (Author, then Year) that was
#    already in effect; data is already in the correct order, this just
#    refreshes the persisted sortState/sortCondition refs.
# ---------------------------------------------------------------------------
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A2:A20"))
$sort.SortFields.Add($ws.Range("E2:E20"))
$sort.SetRange($ws.Range("A2:J20"))
$sort.Apply()

# ---------------------------------------------------------------------------
# 5) Selection / view bookkeeping to match the saved state in the diff.
# ---------------------------------------------------------------------------
$ws.Range("A2:J20").Select()
